# Generate Report for Handoff
# The three "status" rows in each sheet (Overview, zh-cn, de-de) get reshuffled:
#   - ffff7f4f4eab-...md moves into row 2
#   - fffffff9d4d57d-...md moves into row 3
#   - 9705a1c5-...md moves into row 4, and its status flips from
#     "Handed back: in sync with en-US" to "Ready for handoff", along with
#     refreshed handoff timestamps.
# Row 5 (.localization-config) is untouched.

$wb = $excel.ActiveWorkbook

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A3").Value = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A4").Value = "9705a1c5-b117-449f-a566-503268fdd540.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.md" }
}

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md"
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-02-06 04:22:19"
$wsZh.Range("E2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$wsZh.Range("F2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-02-06 04:23:01"
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md"
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-06 04:22:19"
$wsZh.Range("E3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$wsZh.Range("F3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-02-06 04:23:01"
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("A4").Value = "9705a1c5-b117-449f-a566-503268fdd540.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-02-06 04:29:42"
$wsZh.Range("E4").Value = "9705a1c5-b117-449f-a566-503268fdd540.md"
$wsZh.Range("F4").Value = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.zh-cn.xlf"
$wsZh.Range("G4").Value = "2016-02-06 04:28:47"
$wsZh.Range("H4").Value = "Include"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md" }
    elseif ($addr -eq '$C$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf" }
    elseif ($addr -eq '$E$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" }
    elseif ($addr -eq '$F$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.zh-cn.xlf" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.zh-cn.xlf" }
    elseif ($addr -eq '$E$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.md" }
    elseif ($addr -eq '$F$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.zh-cn.xlf" }
}

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md"
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf"
$wsDe.Range("D2").Value = "2016-02-06 04:22:30"
$wsDe.Range("E2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$wsDe.Range("F2").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf"
$wsDe.Range("G2").Value = "2016-02-06 04:23:19"
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md"
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-06 04:22:30"
$wsDe.Range("E3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md"
$wsDe.Range("F3").Value = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf"
$wsDe.Range("G3").Value = "2016-02-06 04:23:19"
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("A4").Value = "9705a1c5-b117-449f-a566-503268fdd540.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.de-de.xlf"
$wsDe.Range("D4").Value = "2016-02-06 04:29:53"
$wsDe.Range("E4").Value = "9705a1c5-b117-449f-a566-503268fdd540.md"
$wsDe.Range("F4").Value = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.de-de.xlf"
$wsDe.Range("G4").Value = "2016-02-06 04:29:06"
$wsDe.Range("H4").Value = "Include"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') { $hl.TextToDisplay = "ffff7f4f4eab-1736-4cf3-bc7e-866986b3265f.md" }
    elseif ($addr -eq '$C$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf" }
    elseif ($addr -eq '$E$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.md" }
    elseif ($addr -eq '$F$2') { $hl.TextToDisplay = "82cd0e9d-47d3-40fb-9dab-c63504f5cbb6.1d1877cdb56ca16a4c696336e6b3429dbb4f778f.de-de.xlf" }
    elseif ($addr -eq '$A$3') { $hl.TextToDisplay = "fffffff9d4d57d-9659-44e2-a62f-74dee04596f3.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.de-de.xlf" }
    elseif ($addr -eq '$E$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.md" }
    elseif ($addr -eq '$F$4') { $hl.TextToDisplay = "9705a1c5-b117-449f-a566-503268fdd540.f4ce14c4a55ed889d8a6ed2a618d432749ba5f19.de-de.xlf" }
}
